# Updates cryptos list data (Price and Volume(1h) columns) to match
# latest scrape, per commit "Updated cryptos list on Sun Oct  6 18:38:51 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.744.94"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "2.441.32"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'566.67"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "'145.69"
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  +2.16%  "
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").Value = "'5.29"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").Value = "'26.91"
$ws.Range("E13").Value = "  +5.08%  "
$ws.Range("E14").Value = "  +5.89%  "
$ws.Range("D15").Value = "2.882.60"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").Value = "62.486.22"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").Value = "2.442.74"
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("D18").Value = "'11.27"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "'6.92"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("D20").Value = "'323.44"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").Value = "'4.17"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "'67.29"
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("D24").Value = "'1.80"
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("E26").Value = "  +8.45%  "
$ws.Range("D27").Value = "'574.39"
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("D28").Value = "2.561.41"
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "'8.41"
$ws.Range("E30").Value = "  +2.69%  "
$ws.Range("E31").Value = "  +3.14%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "'1.87"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("E34").Value = "  +2.25%  "
$ws.Range("E35").Value = "  +3.28%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("D39").Value = "'18.81"
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("D40").Value = "'148.57"
$ws.Range("E40").Value = "  -1.73%  "
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("E43").Value = "  +7.06%  "
$ws.Range("D44").Value = "'147.96"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").Value = "'3.66"
$ws.Range("E45").Value = "  +1.30%  "
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("D47").Value = "'20.49"
$ws.Range("E47").Value = "  +3.10%  "
$ws.Range("E48").Value = "  +2.49%  "
$ws.Range("D49").Value = "'0.0924"
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("E50").Value = "  +2.42%  "
$ws.Range("E51").Value = "  +2.21%  "
